$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 24 (shifts existing rows 24-33 down to 27-36),
# as part of adding the "proceed to checkout" flow rows to the test case sheet.
$ws.Rows("24:26").Insert() | Out-Null

# The inserted rows pick up formatting from the row below by default; re-apply the
# formatting from row 23 (the row above the insertion point) onto the new rows.
$ws.Range("A23:F23").Copy() | Out-Null
$ws.Range("A24:F26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The row that used to be 24 (now row 27, "Apply coupon code or promotional code")
# ends up with a slightly taller custom row height after the edit.
$ws.Rows(27).RowHeight = 16.2

# Update the view: scroll back to the top and select E13.
$ws.Range("E13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
